$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number and report dates) ---
$ws.Range("A8").Characters(21,2).Text = "47"
$ws.Range("C9").Characters(27,10).Text = "11/21/2022"
$ws.Range("C9").Characters(48,10).Text = "11/27/2022"

# --- Cells changing from numeric style to text ("0" / "***.*") style ---
# Donor cell A36 carries the text style (s=14) used for these markers.
$ws.Range("A36").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("A36").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("A36").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("A36").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("A36").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("A36").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("C26").PasteSpecial(-4122)

# --- Cells changing from text style to numeric style ---
# Donor cell C36 carries the numeric style (s=15).
$ws.Range("C36").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 1
$ws.Range("C36").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 1

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -16.666666666666
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -53.571428571428
$ws.Range("I16").Value = 208
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 22.35294117647
$ws.Range("L16").Value = 128.571428571429
$ws.Range("M16").Value = 23.809523809523
$ws.Range("N16").Value = -76.965669988925
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 186
$ws.Range("J17").Value = 181
$ws.Range("K17").Value = 2.762430939226
$ws.Range("L17").Value = 5.681818181818
$ws.Range("M17").Value = 31.914893617021
$ws.Range("N17").Value = -65.619223659889
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 354
$ws.Range("J18").Value = 261
$ws.Range("K18").Value = 35.632183908046
$ws.Range("L18").Value = 47.5
$ws.Range("M18").Value = 76.119402985074
$ws.Range("N18").Value = -58.156028368794
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -29.166666666666
$ws.Range("F19").Value = 82
$ws.Range("G19").Value = 94
$ws.Range("H19").Value = -12.765957446808
$ws.Range("I19").Value = 945
$ws.Range("J19").Value = 723
$ws.Range("K19").Value = 30.705394190871
$ws.Range("L19").Value = 79.657794676806
$ws.Range("M19").Value = 35.38681948424
$ws.Range("N19").Value = -38.636363636363
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 91
$ws.Range("K20").Value = -50.54945054945
$ws.Range("L20").Value = -22.413793103448
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -91.26213592233
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -38.636363636363
$ws.Range("F21").Value = 138
$ws.Range("G21").Value = 168
$ws.Range("H21").Value = -17.857142857142
$ws.Range("I21").Value = 1760
$ws.Range("J21").Value = 1448
$ws.Range("K21").Value = 21.546961325966
$ws.Range("L21").Value = 59.564823209428
$ws.Range("M21").Value = 38.691883372734
$ws.Range("N21").Value = -59.954493742889
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = -33.333333333333
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 120
$ws.Range("J23").Value = 183
$ws.Range("K23").Value = -34.426229508196
$ws.Range("L23").Value = -20.529801324503
$ws.Range("M23").Value = 16.504854368932
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 132
$ws.Range("G24").Value = 151
$ws.Range("H24").Value = -12.582781456953
$ws.Range("I24").Value = 2007
$ws.Range("J24").Value = 1101
$ws.Range("K24").Value = 82.288828337874
$ws.Range("L24").Value = 130.955120828539
$ws.Range("M24").Value = 26.624605678233
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = -7.692307692307
$ws.Range("I25").Value = 429
$ws.Range("J25").Value = 365
$ws.Range("K25").Value = 17.534246575342
$ws.Range("L25").Value = 50
$ws.Range("M25").Value = -0.232558139534
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 32
$ws.Range("J26").Value = 31
$ws.Range("K26").Value = 3.225806451612
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 86
$ws.Range("J27").Value = 65
$ws.Range("K27").Value = 32.307692307692
$ws.Range("L27").Value = 56.363636363636
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = 12.5
$ws.Range("L30").Value = 800
